{"js": "// Update the multiplication-fact answers in the first (and only) table.\n// Each content row (0, 4, 9, 14, 19 in the 20-row grid) holds 5 formula\n// cells \"AAA\u00d7B=CCCC\"; replace every one of the 25 cells with its new\n// value, addressed by absolute row/column index so there is no ambiguity\n// between old/new values that happen to collide (e.g. \"304\u00d77=2128\").\nconst table = context.document.body.tables.getFirst();\n\nconst updates = [\n  { row: 0, col: 0, text: \"592\u00d76=3552\" },\n  { row: 0, col: 1, text: \"175\u00d72=350\" },\n  { row: 0, col: 2, text: \"313\u00d75=1565\" },\n  { row: 0, col: 3, text: \"342\u00d76=2052\" },\n  { row: 0, col: 4, text: \"579\u00d74=2316\" },\n\n  { row: 4, col: 0, text: \"558\u00d77=3906\" },\n  { row: 4, col: 1, text: \"998\u00d72=1996\" },\n  { row: 4, col: 2, text: \"667\u00d72=1334\" },\n  { row: 4, col: 3, text: \"877\u00d78=7016\" },\n  { row: 4, col: 4, text: \"716\u00d76=4296\" },\n\n  { row: 9, col: 0, text: \"536\u00d74=2144\" },\n  { row: 9, col: 1, text: \"285\u00d73=855\" },\n  { row: 9, col: 2, text: \"997\u00d73=2991\" },\n  { row: 9, col: 3, text: \"453\u00d73=1359\" },\n  { row: 9, col: 4, text: \"663\u00d72=1326\" },\n\n  { row: 14, col: 0, text: \"673\u00d74=2692\" },\n  { row: 14, col: 1, text: \"323\u00d79=2907\" },\n  { row: 14, col: 2, text: \"902\u00d73=2706\" },\n  { row: 14, col: 3, text: \"455\u00d77=3185\" },\n  { row: 14, col: 4, text: \"443\u00d78=3544\" },\n\n  { row: 19, col: 0, text: \"692\u00d73=2076\" },\n  { row: 19, col: 1, text: \"304\u00d77=2128\" },\n  { row: 19, col: 2, text: \"821\u00d79=7389\" },\n  { row: 19, col: 3, text: \"732\u00d73=2196\" },\n  { row: 19, col: 4, text: \"961\u00d79=8649\" },\n];\n\nfor (const u of updates) {\n  table.getCell(u.row, u.col).value = u.text;\n}\n\nawait context.sync();\n", "ps1": "# Update the multiplication-fact answers in the first (and only) table.\n# Each content row (1, 5, 10, 15, 20 in the 1-indexed 20-row grid) holds 5\n# formula cells \"AAA\u00d7B=CCCC\"; replace every one of the 25 cells with its\n# new value, addressed by absolute 1-based row/column index so there is no\n# ambiguity between old/new values that happen to collide (e.g. \"304\u00d77=2128\").\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$t.Cell(1,1).Range.Text = \"592\u00d76=3552\"\n$t.Cell(1,2).Range.Text = \"175\u00d72=350\"\n$t.Cell(1,3).Range.Text = \"313\u00d75=1565\"\n$t.Cell(1,4).Range.Text = \"342\u00d76=2052\"\n$t.Cell(1,5).Range.Text = \"579\u00d74=2316\"\n\n$t.Cell(5,1).Range.Text = \"558\u00d77=3906\"\n$t.Cell(5,2).Range.Text = \"998\u00d72=1996\"\n$t.Cell(5,3).Range.Text = \"667\u00d72=1334\"\n$t.Cell(5,4).Range.Text = \"877\u00d78=7016\"\n$t.Cell(5,5).Range.Text = \"716\u00d76=4296\"\n\n$t.Cell(10,1).Range.Text = \"536\u00d74=2144\"\n$t.Cell(10,2).Range.Text = \"285\u00d73=855\"\n$t.Cell(10,3).Range.Text = \"997\u00d73=2991\"\n$t.Cell(10,4).Range.Text = \"453\u00d73=1359\"\n$t.Cell(10,5).Range.Text = \"663\u00d72=1326\"\n\n$t.Cell(15,1).Range.Text = \"673\u00d74=2692\"\n$t.Cell(15,2).Range.Text = \"323\u00d79=2907\"\n$t.Cell(15,3).Range.Text = \"902\u00d73=2706\"\n$t.Cell(15,4).Range.Text = \"455\u00d77=3185\"\n$t.Cell(15,5).Range.Text = \"443\u00d78=3544\"\n\n$t.Cell(20,1).Range.Text = \"692\u00d73=2076\"\n$t.Cell(20,2).Range.Text = \"304\u00d77=2128\"\n$t.Cell(20,3).Range.Text = \"821\u00d79=7389\"\n$t.Cell(20,4).Range.Text = \"732\u00d73=2196\"\n$t.Cell(20,5).Range.Text = \"961\u00d79=8649\"\n"}
